$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LA PRESSE")

# Column F builds a JSON fragment from columns A-D, e.g.:
#   ="{ ""foreign"": """ & A2 & """, ""grammar"": """ & B2 & """, ""pronunciation"": """ & C2 & """, ""meaning"": """ & D2 & """ },"
# Row 2 already has this formula. Extend it down through row 17 (new vocabulary rows).
for ($r = 3; $r -le 17; $r++) {
    $formula = '= "{ ""foreign"": """ & A' + $r + ' & """, ""grammar"": """ & B' + $r + ' & """, ""pronunciation"": """ & C' + $r + ' & """, ""meaning"": """ & D' + $r + ' & """ },"'
    $ws.Range("F$r").Formula = $formula
}

# Leave the selection over the freshly filled formula range, as happens after
# filling a formula down a column in Excel.
$ws.Range("F3:F17").Select()
